$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the previous month's block (rows 11-19)
# down onto the new rows 20-28 so the new cells inherit the same number
# format (date column C) and font style (column D) as the existing data.
$ws.Range("A11:D19").Copy($ws.Range("A20:D28"))

# New charge entries for the September 2024 (serial 45536) period, in the
# same row order used by the previous monthly blocks.
$rows = @(
    @("electricite", 71),
    @("gaz", 22),
    @("copro", 281),
    @("box ", 30),
    @("credit et assurances", 400),
    @("marceline", 0),
    @("charbel", 0),
    @("adekemi", 0),
    @("kadi", 0)
)

$r = 20
foreach ($item in $rows) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = 45536
    $ws.Cells.Item($r, 4).Value = "Appartement Madoumier T4"
    $r = $r + 1
}

# Update the active selection/view to match the saved workbook state.
$ws.Range("I30").Select()
